$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Product Price" cell (D2) from the text value "31500"
# to the proper numeric value 350.
$ws.Range("D2").Value = 350
